# Updated cryptos list on Thu Apr  4 23:26:02 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures scraped
# from coinranking.com for each coin row, and reflects that Mantle (MNT)
# overtook Stellar (XLM) in the ranking, so rows 47/48 swap contents.
#
# Prices in column D are stored as plain text (not numbers) in the workbook,
# so for any value that Excel's automatic type inference would otherwise
# coerce into a floating point number (and thereby mangle through binary
# rounding, e.g. 583.27 -> 583.26999999999998), the cell's number format is
# first forced to Text ("@") before the value is written, guaranteeing the
# text is stored byte-for-byte.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.229.44"
$ws.Range("E2").Value = "  +3.07%  "
$ws.Range("D3").Value = "3.322.28"
$ws.Range("E3").Value = "  +0.04%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.27"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "183.17"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +2.64%  "
$ws.Range("D9").Value = "3.319.37"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.37"
$ws.Range("E13").Value = "  +3.52%  "
$ws.Range("D14").Value = "635.75"
$ws.Range("E14").Value = "  +6.85%  "
$ws.Range("D15").Value = "3.858.23"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("D17").Value = "68.325.38"
$ws.Range("E17").Value = "  +3.23%  "
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "3.324.90"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("D21").Value = "10.94"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("D23").Value = "17.66"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.12"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("E26").Value = "  +0.35%  "
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.59"
$ws.Range("E28").Value = "  +1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.51"
$ws.Range("E29").Value = "  +5.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.60"
$ws.Range("E30").Value = "  +1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.70"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "607.68"
$ws.Range("E32").Value = "  +8.65%  "
$ws.Range("D33").Value = "3.941.20"
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.97"
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("D36").Value = "3.52"
$ws.Range("E36").Value = "  -4.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.83"
$ws.Range("E38").Value = "  -0.08%  "
$ws.Range("D39").Value = "3.27"
$ws.Range("E39").Value = "  +4.12%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("E41").Value = "  +3.37%  "
$ws.Range("D42").Value = "32.68"
$ws.Range("E42").Value = "  -1.77%  "
$ws.Range("D43").Value = "0.0₃0687"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.338"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.42"
$ws.Range("E47").Value = "  +14.73%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.129"
$ws.Range("E48").Value = "  +1.89%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.56"
$ws.Range("E50").Value = "  +1.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "131.24"
$ws.Range("E51").Value = "  +1.29%  "
